# Inserts two new price rows (Sandia - "Extra" and "Primera" quality,
# dated serial 44551 / 2021-12-21) above the existing row 93, pushing the
# old rows 93-110 down to 95-112 (matches new dimension A1:R112).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 93 (shifts 93-110 -> 95-112).
$ws.Rows.Item(93).Insert()
$ws.Rows.Item(93).Insert()

# New row 93: Sandia, "Extra" quality.
$ws.Range("A93").Value = 5
$ws.Range("B93").Value = "Macroferia Regional de Talca"
$ws.Range("C93").Value = "Maule"
$ws.Range("D93").Value = 44551
$ws.Range("E93").Value = 7
$ws.Range("F93").Value = 100112028
$ws.Range("G93").Value = "Sandia"
$ws.Range("H93").Value = "Sin especificar"
$ws.Range("I93").Value = "Extra"
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 3000
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = 3000
$ws.Range("N93").Value = "$/unidad"
$ws.Range("O93").Value = "Región del Maule"
$ws.Range("P93").Value = 3000
$ws.Range("Q93").Value = 1
$ws.Range("R93").Value = "Hortaliza"

# New row 94: Sandia, "Primera" quality.
$ws.Range("A94").Value = 5
$ws.Range("B94").Value = "Macroferia Regional de Talca"
$ws.Range("C94").Value = "Maule"
$ws.Range("D94").Value = 44551
$ws.Range("E94").Value = 7
$ws.Range("F94").Value = 100112028
$ws.Range("G94").Value = "Sandia"
$ws.Range("H94").Value = "Sin especificar"
$ws.Range("I94").Value = "Primera"
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 2500
$ws.Range("L94").Value = 2500
$ws.Range("M94").Value = 2500
$ws.Range("N94").Value = "$/unidad"
$ws.Range("O94").Value = "Región del Maule"
$ws.Range("P94").Value = 2500
$ws.Range("Q94").Value = 1
$ws.Range("R94").Value = "Hortaliza"
